# Recalculated market-data snapshot refresh (scheduled runner).
# Updates cached price/profit figures (columns H-N) across several
# leve-profit worksheets to reflect newer Market Board averages.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 861.6667
$ws.Range("J32").Value = 861.6667
$ws.Range("L32").Value = 861.6667
$ws.Range("N32").Value = -1513.6667
$ws.Range("H107").Value = 991
$ws.Range("I107").Value = 990.3077
$ws.Range("K107").Value = 990.3077
$ws.Range("M107").Value = 929.6923
$ws.Range("H112").Value = 1212.8772
$ws.Range("J112").Value = 1212.8772
$ws.Range("L112").Value = 3638.6316
$ws.Range("N112").Value = -5854.6316
$ws.Range("H132").Value = 4103661.8
$ws.Range("I132").Value = 4722422
$ws.Range("K132").Value = 14167266
$ws.Range("M132").Value = -14164736
$ws.Range("H138").Value = 3067.9048
$ws.Range("I138").Value = 1861.3182
$ws.Range("J138").Value = 3715.3416
$ws.Range("K138").Value = 5583.9546
$ws.Range("L138").Value = 11146.0248
$ws.Range("M138").Value = -443.9546
$ws.Range("N138").Value = -21426.0248
$ws.Range("H141").Value = 2540.4783
$ws.Range("J141").Value = 2843.5715
$ws.Range("L141").Value = 8530.7145
$ws.Range("N141").Value = -18890.7145

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 180000
$ws.Range("J11").Value = 180000
$ws.Range("L11").Value = 180000
$ws.Range("N11").Value = -180288
$ws.Range("H45").Value = 334550
$ws.Range("I45").Value = 334550
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 334550
$ws.Range("L45").Value = 0
$ws.Range("N45").Value = -334173
$ws.Range("H61").Value = 1239
$ws.Range("I61").Value = 1239
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1239
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = -1027
$ws.Range("H74").Value = 963.5789
$ws.Range("I74").Value = 892
$ws.Range("K74").Value = 892
$ws.Range("M74").Value = -18
$ws.Range("H77").Value = 963.5789
$ws.Range("I77").Value = 892
$ws.Range("K77").Value = 4460
$ws.Range("M77").Value = -92
$ws.Range("H132").Value = 14330.805
$ws.Range("I132").Value = 15541.951
$ws.Range("K132").Value = 46625.853
$ws.Range("M132").Value = -44095.853
$ws.Range("H135").Value = 41801.5
$ws.Range("J135").Value = 41801.5
$ws.Range("L135").Value = 41801.5
$ws.Range("N135").Value = -51941.5
$ws.Range("H136").Value = 1239
$ws.Range("I136").Value = 1239
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3717
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = -1167
$ws.Range("M45").ClearContents()
$ws.Range("M61").ClearContents()
$ws.Range("M136").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 47641610
$ws.Range("I107").Value = 50023616
$ws.Range("K107").Value = 50023616
$ws.Range("M107").Value = -50021696
$ws.Range("H134").Value = 1839.7683
$ws.Range("I134").Value = 1704.5135
$ws.Range("K134").Value = 5113.5405
$ws.Range("M134").Value = -2578.5405

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 56326.055
$ws.Range("I16").Value = 795.0909
$ws.Range("J16").Value = 143589
$ws.Range("K16").Value = 795.0909
$ws.Range("L16").Value = 143589
$ws.Range("M16").Value = -508.0909
$ws.Range("N16").Value = -144163
$ws.Range("H62").Value = 2460
$ws.Range("I62").Value = 2242.8572
$ws.Range("J62").Value = 2650
$ws.Range("K62").Value = 2242.8572
$ws.Range("L62").Value = 2650
$ws.Range("M62").Value = -1618.8572
$ws.Range("N62").Value = -3898
$ws.Range("H65").Value = 2460
$ws.Range("I65").Value = 2242.8572
$ws.Range("J65").Value = 2650
$ws.Range("K65").Value = 11214.286
$ws.Range("L65").Value = 13250
$ws.Range("M65").Value = -8094.286
$ws.Range("N65").Value = -19490
$ws.Range("H113").Value = 56326.055
$ws.Range("I113").Value = 795.0909
$ws.Range("J113").Value = 143589
$ws.Range("K113").Value = 795.0909
$ws.Range("L113").Value = 143589
$ws.Range("M113").Value = 1374.9091
$ws.Range("N113").Value = -147929

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 127337.75
$ws.Range("I70").Value = 251700.5
$ws.Range("K70").Value = 755101.5
$ws.Range("M70").Value = -754786.5
$ws.Range("H73").Value = 127337.75
$ws.Range("I73").Value = 251700.5
$ws.Range("K73").Value = 755101.5
$ws.Range("M73").Value = -754009.5
$ws.Range("H80").Value = 8641.857
$ws.Range("J80").Value = 9229.691999999999
$ws.Range("L80").Value = 27689.076
$ws.Range("N80").Value = -29561.076
$ws.Range("H83").Value = 8641.857
$ws.Range("J83").Value = 9229.691999999999
$ws.Range("L83").Value = 83067.22799999999
$ws.Range("N83").Value = -92427.22799999999
$ws.Range("H127").Value = 1047.5834
$ws.Range("J127").Value = 1066.909
$ws.Range("L127").Value = 3200.727
$ws.Range("N127").Value = -13120.727
$ws.Range("H131").Value = 791689.4
$ws.Range("J131").Value = 955871.9
$ws.Range("L131").Value = 2867615.7
$ws.Range("N131").Value = -2877695.7

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 594651.4399999999
$ws.Range("I107").Value = 372.2
$ws.Range("J107").Value = 1443621.8
$ws.Range("K107").Value = 372.2
$ws.Range("L107").Value = 1443621.8
$ws.Range("M107").Value = 1547.8
$ws.Range("N107").Value = -1447461.8
$ws.Range("H132").Value = 4156.8184
$ws.Range("I132").Value = 3123.6667
$ws.Range("K132").Value = 9371.000100000001
$ws.Range("M132").Value = -6841.000100000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6112.857
$ws.Range("I46").Value = 690
$ws.Range("J46").Value = 7016.6665
$ws.Range("K46").Value = 690
$ws.Range("L46").Value = 7016.6665
$ws.Range("M46").Value = -502
$ws.Range("N46").Value = -7392.6665
$ws.Range("H100").Value = 1771.4286
$ws.Range("I100").Value = 1720
$ws.Range("J100").Value = 1900
$ws.Range("K100").Value = 1720
$ws.Range("L100").Value = 1900
$ws.Range("M100").Value = -1179
$ws.Range("N100").Value = -2982
$ws.Range("H132").Value = 8412.25
$ws.Range("I132").Value = 9216.666999999999
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 27650.001
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -25120.001
$ws.Range("N132").Value = -23057

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 470.54166
$ws.Range("I113").Value = 350.27777
$ws.Range("J113").Value = 831.3333
$ws.Range("K113").Value = 1050.83331
$ws.Range("L113").Value = 2493.9999
$ws.Range("M113").Value = 1119.16669
$ws.Range("N113").Value = -6833.9999
$ws.Range("H132").Value = 15134.818
$ws.Range("I132").Value = 9936
$ws.Range("K132").Value = 29808
$ws.Range("M132").Value = -27278
$ws.Range("H136").Value = 12880.021
$ws.Range("I136").Value = 19346.264
$ws.Range("J136").Value = 4521.2197
$ws.Range("K136").Value = 58038.792
$ws.Range("L136").Value = 13563.6591
$ws.Range("M136").Value = -55488.792
$ws.Range("N136").Value = -18663.6591
